$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C6 entirely
$ws.Range("C6").ClearContents()

# Update B6 value and style
$ws.Range("B6").Value = "Cent. vyrovnanie"
$ws.Range("B6").Style = $ws.Range("B5").Style

# Move selection to B11
$ws.Range("B11").Select()
